$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.491.71'
$ws.Range("E2").Value = '  -1.05%  '
$ws.Range("D3").Value = '2.929.59'
$ws.Range("E3").Value = '  -2.54%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '''374.52'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.70%  '
$ws.Range("D6").Value = '''103.45'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.22%  '
$ws.Range("E7").Value = '  -2.89%  '
$ws.Range("E8").Value = '  -0.20%  '
$ws.Range("D9").Value = '''0.586'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.16%  '
$ws.Range("D10").Value = '''36.92'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.97%  '
$ws.Range("E11").Value = '  -0.58%  '
$ws.Range("D12").Value = '''0.0838'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.38%  '
$ws.Range("E13").Value = '  -3.50%  '
$ws.Range("D14").Value = '3.390.64'
$ws.Range("E14").Value = '  -2.57%  '
$ws.Range("E15").Value = '  -3.56%  '
$ws.Range("D16").Value = '2.923.77'
$ws.Range("E16").Value = '  -2.36%  '
$ws.Range("D17").Value = '''0.928'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -8.44%  '
$ws.Range("D18").Value = '51.459.24'
$ws.Range("D19").Value = '''3.42'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.58%  '
$ws.Range("D20").Value = '''7.32'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.91%  '
$ws.Range("D21").Value = '''12.95'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.41%  '
$ws.Range("D22").Value = '0.0₃0945'
$ws.Range("E22").Value = '  -2.68%  '
$ws.Range("D23").Value = '''68.28'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.22%  '
$ws.Range("D24").Value = '''262.13'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.60%  '
$ws.Range("E25").Value = '  +1.01%  '
$ws.Range("E26").Value = '  -5.65%  '
$ws.Range("E27").Value = '  -5.10%  '
$ws.Range("E28").Value = '  +0.00%  '
$ws.Range("D29").Value = '''25.75'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.57%  '
$ws.Range("D30").Value = '''7.29'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.16%  '
$ws.Range("D31").Value = '''6.91'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.56%  '
$ws.Range("E32").Value = '  -5.24%  '
$ws.Range("D33").Value = '''9.80'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.76%  '
$ws.Range("E34").Value = '  -3.30%  '
$ws.Range("D35").Value = '''51.05'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.16%  '
$ws.Range("D36").Value = '''33.97'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.77%  '
$ws.Range("E37").Value = '  +0.37%  '
$ws.Range("D38").Value = '''0.0423'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.04%  '
$ws.Range("E39").Value = '  -9.50%  '
$ws.Range("D40").Value = '''16.96'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.72%  '
$ws.Range("D41").Value = '''2.58'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -9.50%  '
$ws.Range("D42").Value = '''1.81'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.31%  '
$ws.Range("E43").Value = '  -2.41%  '
$ws.Range("D44").Value = '''123.64'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.48%  '
$ws.Range("D45").Value = '''21.76'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.37%  '
$ws.Range("E46").Value = '  -4.03%  '
$ws.Range("D47").Value = '''0.271'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +11.11%  '
$ws.Range("D48").Value = '2.019.75'
$ws.Range("E48").Value = '  -4.83%  '
$ws.Range("E49").Value = '  -1.72%  '
$ws.Range("E50").Value = '  -5.22%  '
$ws.Range("D51").Value = '3.210.11'
$ws.Range("E51").Value = '  -2.86%  '
